# The sheet's only real data edit: the "LTS-Tollgrade" row's Date (B11)
# moves from 2024-07-13 to 2023-07-13 (year corrected 2024 -> 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "2023-07-13"

# View-state changes captured in the diff: the sheet was scrolled so row 3
# is the top visible row, and the active selection moved from D6 to D10.
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
